$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avichal")
$ws.Range("A2").Value = "test"
